$d = $word.ActiveDocument

$titleText  = "Play Age of the Gods: God of Storms for Free | Review"
$metaSuffix = ": Get an unbiased review of Age of the Gods: God of Storms slot game and play for free. Learn about the gameplay, graphics, jackpots, and more."
$metaBold   = "Meta description"
$oldItalicText = "Get an unbiased review of Age of the Gods: God of Storms slot game and play for free. Learn about the gameplay, graphics, jackpots, and more."
$newItalicText = 'Create a feature image for "Age of the Gods God of Storms": - Draw a cartoon-style image featuring a happy Maya warrior with glasses. - The warrior should stand triumphantly next to the game''s title. - The background should be a stormy sea with waves crashing in the distance. - The warrior should hold a sword and shield, ready to take on the storms and win big.'

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph (the very first paragraph of the document). Built via
#    InsertXML so the resulting run layout matches exactly: a leading
#    empty run, a bold "Meta description" run, and a plain run with the
#    rest of the sentence.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
if ($titleRange.Text.TrimEnd([char]13, [char]7) -ne $titleText) {
    throw "Unexpected first paragraph text: $($titleRange.Text)"
}
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $metaBold + '</w:t></w:r><w:r><w:t>' + $metaSuffix + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicate bold title
#    paragraph entirely, and replace the italic meta-description
#    paragraph's text with the new image prompt (keeping the leading
#    empty run + italic run structure).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
if ($dupTitlePara.Range.Text.TrimEnd([char]13, [char]7) -ne $titleText) {
    throw "Unexpected second-to-last paragraph text: $($dupTitlePara.Range.Text)"
}
$dupTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$imgPromptPara = $d.Paragraphs($count)
if ($imgPromptPara.Range.Text.TrimEnd([char]13, [char]7) -ne $oldItalicText) {
    throw "Unexpected last paragraph text: $($imgPromptPara.Range.Text)"
}
$imgPromptPara.Range.Delete()

$imgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newItalicText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs($count).Range.InsertXML($imgXml)

Write-Output "done"
